# Automatic update of files.
# Applies the taxon-sort-order bump for Garnlav (taxon 6425) 79243 -> 79244,
# the analogous bump on row 15, and restores the correct ordering of the two
# swapped-row pairs (10/11 and 21/22) that the source export shuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple Taxonsorteringsordning (column B) bumps: 79243 -> 79244
# ---------------------------------------------------------------------
foreach ($r in 5, 6, 7, 8) {
    $ws.Cells.Item($r, 2).Value = 79244
}

# ---------------------------------------------------------------------
# 2) Swap the full data rows 10 and 11 (every column A..AY), since the
#    export had them in the wrong order.
# ---------------------------------------------------------------------
function Swap-Rows($ws, $rowA, $rowB, $maxCol) {
    $valsA = @()
    $valsB = @()

    for ($c = 1; $c -le $maxCol; $c++) {
        $valsA += ,$ws.Cells.Item($rowA, $c).Value2
        $valsB += ,$ws.Cells.Item($rowB, $c).Value2
    }

    # Columns Y (25) and AA (27) hold dates stored as plain text
    # ("2026-02-14"); force text format first so the COM Value setter
    # doesn't "helpfully" reinterpret them as real date serials.
    foreach ($dateCol in 25, 27) {
        $ws.Cells.Item($rowA, $dateCol).NumberFormat = "@"
        $ws.Cells.Item($rowB, $dateCol).NumberFormat = "@"
    }

    for ($c = 1; $c -le $maxCol; $c++) {
        $ws.Cells.Item($rowA, $c).Value = $valsB[$c - 1]
        $ws.Cells.Item($rowB, $c).Value = $valsA[$c - 1]
    }
}

Swap-Rows $ws 10 11 51

# Row 10 additionally needs the Taxonsorteringsordning bump (the value that
# lands on row 10 is the Garnlav one, 79243 -> 79244), matching change (1).
$ws.Cells.Item(10, 2).Value = 79244

# ---------------------------------------------------------------------
# 3) Swap the full data rows 21 and 22 (every column A..AY) similarly.
# ---------------------------------------------------------------------
Swap-Rows $ws 21 22 51

# ---------------------------------------------------------------------
# 4) Simple Taxonsorteringsordning (column B) bump on row 15.
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 2).Value = 91829
